$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion rates text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $ws1.Range("A1")
$text = $cellA1.Value()
$text = $text.Replace("1000 Bs = 6.02 = 24855.54 pesos", "1000 Bs = 6.04 = 24821.81 pesos")
$text = $text.Replace("24855.54 pesos = 6.02 = 967.21 Bs", "24821.81 pesos = 6.02 = 976.01 Bs")
$cellA1.Value = $text

# --- Sheet "tasas": update the rate cells N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 165.5
$ws2.Range("O10").Value = 4108.01
$ws2.Range("N12").Value = 4119.99
$ws2.Range("O12").Value = 162
